$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.615.92"
$ws.Range("E2").Value = "  -1.29%  "

$ws.Range("D3").Value = "1.849.28"
$ws.Range("E3").Value = "  -2.39%  "

$ws.Range("E4").Value = "  -0.64%  "

$ws.Range("D5").Formula = "=""333.79"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").Formula = "=""1.010"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -0.65%  "

$ws.Range("D7").Formula = "=""0.4592"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -2.22%  "

$ws.Range("D8").Formula = "=""0.3894"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("D9").Formula = "=""45.98"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -3.62%  "

$ws.Range("D10").Formula = "=""0.07916"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -1.70%  "

$ws.Range("D11").Formula = "=""1.001"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -1.99%  "

$ws.Range("D12").Formula = "=""21.58"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -1.18%  "

$ws.Range("D13").Value = "1.857.17"
$ws.Range("E13").Value = "  -2.31%  "

$ws.Range("D14").Formula = "=""5.953"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").Formula = "=""7.182"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("D16").Formula = "=""1.012"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -0.75%  "

$ws.Range("D17").Formula = "=""88.52"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +1.07%  "

$ws.Range("D18").Formula = "=""0.06716"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").Formula = "=""0.00001037"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -1.80%  "

$ws.Range("D20").Formula = "=""17.16"""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -0.21%  "

$ws.Range("D21").Formula = "=""1.011"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -0.48%  "

$ws.Range("D22").Value = "27.600.17"
$ws.Range("E22").Value = "  -1.42%  "

$ws.Range("D23").Formula = "=""5.414"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -1.92%  "

$ws.Range("D24").Formula = "=""10.90"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.84%  "

$ws.Range("E25").Value = "  -1.97%  "

$ws.Range("D26").Formula = "=""159.36"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.06%  "

$ws.Range("D27").Formula = "=""19.55"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -2.57%  "

$ws.Range("D28").Formula = "=""2.126"""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +1.97%  "

$ws.Range("D29").Formula = "=""5.431"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -0.61%  "

$ws.Range("D30").Formula = "=""121.56"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("D31").Formula = "=""0.9730"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +0.19%  "

$ws.Range("D32").Formula = "=""0.09392"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -1.33%  "

$ws.Range("D33").Formula = "=""3.620"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -1.49%  "

$ws.Range("D34").Formula = "=""5.299"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -1.36%  "

$ws.Range("D35").Formula = "=""1.333"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -4.81%  "

$ws.Range("D36").Formula = "=""0.02226"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -1.42%  "

$ws.Range("D37").Formula = "=""0.05996"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -2.36%  "

$ws.Range("D38").Formula = "=""8.354"""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +2.40%  "

$ws.Range("D39").Formula = "=""1.190"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -2.49%  "

$ws.Range("D40").Formula = "=""1.010"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -0.47%  "

$ws.Range("D41").Formula = "=""0.5903"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -1.66%  "

$ws.Range("D42").Formula = "=""10.42"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("D43").Formula = "=""0.1860"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.69%  "

$ws.Range("D44").Formula = "=""1.244"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -2.05%  "

$ws.Range("D45").Formula = "=""0.5570"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -2.43%  "

$ws.Range("D46").Formula = "=""12.09"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -1.53%  "

$ws.Range("D47").Formula = "=""1.910"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -1.61%  "

$ws.Range("D48").Formula = "=""0.06703"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -3.40%  "

$ws.Range("D49").Formula = "=""111.06"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -2.69%  "

$ws.Range("D50").Formula = "=""1.049"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -2.05%  "

$ws.Range("D51").Formula = "=""1.011"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -0.55%  "

Write-Output "done"